# Apply "model and template with unit, description and enum" change:
#  - Row 1 (field names): H1/I1 swap so ExtractionType comes before SamplePortion.
#  - Row 2 (type annotations): now carries units for the float columns; the
#    analyte columns (J:AX) use the new "unit:µmole/ml" annotation and the
#    SamplePortion column (I) uses "unit:µlormg"; the ExtractionType column
#    (H) switches from float to string, matching its swapped header.
#  - Row 3 (new): French field descriptions for the first 9 metadata
#    columns (A:I); the analyte columns (J:AX) get an empty placeholder
#    string so the shared-string table keeps one entry per used cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header names -----------------------------------------------
# Only H1/I1 actually change value (ExtractionType now precedes
# SamplePortion); everything else in row 1 stays as-is.
$ws.Range("H1").Value = "ExtractionType"
$ws.Range("I1").Value = "SamplePortion"

# --- Row 2: type annotations --------------------------------------------
# A:G keep "#string"/"#date" as before; H switches to "#string" (it is now
# the ExtractionType column) and I becomes the SamplePortion float w/ unit.
$ws.Range("A2").Value = "#string"
$ws.Range("B2").Value = "#string"
$ws.Range("C2").Value = "#date"
$ws.Range("D2").Value = "#string"
$ws.Range("E2").Value = "#string"
$ws.Range("F2").Value = "#string"
$ws.Range("G2").Value = "#string"
$ws.Range("H2").Value = "#string"
$ws.Range("I2").Value = "#float,  unit:µlormg"

# J2:AX2 (the amino-acid concentration columns) all take the molar-unit
# float annotation.
for ($c = 10; $c -le 50; $c++) {
    $ws.Cells.Item(2, $c).Value = "#float,  unit:µmole/ml"
}

# --- Row 3: new description row -----------------------------------------
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#TypeExtraction"
$ws.Range("I3").Value = "#PriseEssai"

# J3:AX3 get an empty string placeholder (matches the blank <t/> shared
# string entry added for these cells).
for ($c = 10; $c -le 50; $c++) {
    $ws.Cells.Item(3, $c).Value = ""
}
